$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 4491.067
$ws.Range("I74").Value = 4209.091
$ws.Range("J74").Value = 5266.5
$ws.Range("K74").Value = 4209.091
$ws.Range("L74").Value = 5266.5
$ws.Range("M74").Value = -3273.091
$ws.Range("N74").Value = -7138.5
# Row 77
$ws.Range("H77").Value = 4491.067
$ws.Range("I77").Value = 4209.091
$ws.Range("J77").Value = 5266.5
$ws.Range("K77").Value = 21045.455
$ws.Range("L77").Value = 26332.5
$ws.Range("M77").Value = -16365.455
$ws.Range("N77").Value = -35692.5
# Row 137
$ws.Range("H137").Value = 1445.7188
$ws.Range("I137").Value = 1326.44
$ws.Range("J137").Value = 1871.7142
$ws.Range("K137").Value = 3979.32
$ws.Range("L137").Value = 5615.142599999999
$ws.Range("M137").Value = -1429.32
$ws.Range("N137").Value = -10715.1426

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
# Row 132
$ws.Range("H132").Value = 3079.8057
$ws.Range("I132").Value = 1959.4348
$ws.Range("J132").Value = 5062
$ws.Range("K132").Value = 5878.3044
$ws.Range("L132").Value = 15186
$ws.Range("M132").Value = -3348.3044
$ws.Range("N132").Value = -20246

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1864.2759
$ws.Range("I94").Value = 1534.2
$ws.Range("K94").Value = 1534.2
$ws.Range("M94").Value = -1083.2

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 30250.25
$ws.Range("J4").Value = 34000
$ws.Range("L4").Value = 34000
$ws.Range("N4").Value = -34224
# Row 22
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -650
$ws.Range("N22").Value = -1700
# Row 31
$ws.Range("H31").Value = 14289907
$ws.Range("I31").Value = 2303.647
$ws.Range("J31").Value = 27783754
$ws.Range("K31").Value = 2303.647
$ws.Range("L31").Value = 27783754
$ws.Range("M31").Value = -2008.647
$ws.Range("N31").Value = -27784344
# Row 34
$ws.Range("H34").Value = 14289907
$ws.Range("I34").Value = 2303.647
$ws.Range("J34").Value = 27783754
$ws.Range("K34").Value = 2303.647
$ws.Range("L34").Value = 27783754
$ws.Range("M34").Value = -2101.647
$ws.Range("N34").Value = -27784158
# Row 105
$ws.Range("H105").Value = 15874599
$ws.Range("I105").Value = 22223872
$ws.Range("J105").Value = 1416.6666
$ws.Range("K105").Value = 22223872
$ws.Range("L105").Value = 1416.6666
$ws.Range("M105").Value = -22222125
$ws.Range("N105").Value = -4910.6666
# Row 134
$ws.Range("H134").Value = 11112498
$ws.Range("I134").Value = 13334681
$ws.Range("J134").Value = 1582.8
$ws.Range("K134").Value = 40004043
$ws.Range("L134").Value = 4748.4
$ws.Range("M134").Value = -40001508
$ws.Range("N134").Value = -9818.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 7645.8237
$ws.Range("I3").Value = 10089.857
$ws.Range("J3").Value = 5935
$ws.Range("K3").Value = 30269.571
$ws.Range("L3").Value = 17805
$ws.Range("M3").Value = -30157.571
$ws.Range("N3").Value = -18029
# Row 23
$ws.Range("H23").Value = 6250093.5
$ws.Range("J23").Value = 97.454544
$ws.Range("L23").Value = 292.363632
$ws.Range("N23").Value = -762.3636320000001
# Row 32
$ws.Range("H32").Value = 3533.3333
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 3850
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 11550
$ws.Range("M32").Value = -2717
$ws.Range("N32").Value = -12116
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
# Row 41
$ws.Range("H41").Value = 2242.8572
$ws.Range("I41").Value = 1300
$ws.Range("J41").Value = 2400
$ws.Range("K41").Value = 3900
$ws.Range("L41").Value = 7200
$ws.Range("M41").Value = -3562
$ws.Range("N41").Value = -7876
# Row 58
$ws.Range("H58").Value = 3000.353
$ws.Range("I58").Value = 1300
$ws.Range("J58").Value = 3106.625
$ws.Range("K58").Value = 3900
$ws.Range("L58").Value = 9319.875
$ws.Range("M58").Value = -3772
$ws.Range("N58").Value = -9575.875
# Row 61
$ws.Range("H61").Value = 441.66666
$ws.Range("J61").Value = 518
$ws.Range("L61").Value = 1554
$ws.Range("N61").Value = -1984
# Row 62
$ws.Range("H62").Value = 6498.3335
$ws.Range("J62").Value = 6498.3335
$ws.Range("L62").Value = 19495.0005
$ws.Range("N62").Value = -20867.0005
# Row 63
$ws.Range("H63").Value = 799
$ws.Range("I63").Value = 799
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2397
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1648
$ws.Range("N63").ClearContents()
# Row 65
$ws.Range("H65").Value = 6498.3335
$ws.Range("J65").Value = 6498.3335
$ws.Range("L65").Value = 58485.0015
$ws.Range("N65").Value = -65349.0015
# Row 66
$ws.Range("H66").Value = 799
$ws.Range("I66").Value = 799
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 7191
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -3447
$ws.Range("N66").ClearContents()
# Row 70
$ws.Range("H70").Value = 2093.182
$ws.Range("I70").Value = 860.7143
$ws.Range("J70").Value = 4250
$ws.Range("K70").Value = 2582.1429
$ws.Range("L70").Value = 12750
$ws.Range("M70").Value = -2267.1429
$ws.Range("N70").Value = -13380
# Row 73
$ws.Range("H73").Value = 2093.182
$ws.Range("I73").Value = 860.7143
$ws.Range("J73").Value = 4250
$ws.Range("K73").Value = 2582.1429
$ws.Range("L73").Value = 12750
$ws.Range("M73").Value = -1490.1429
$ws.Range("N73").Value = -14934
# Row 74
$ws.Range("H74").Value = 20000
$ws.Range("J74").Value = 20000
$ws.Range("L74").Value = 60000
$ws.Range("N74").Value = -62122
# Row 75
$ws.Range("H75").Value = 142857140
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 142857140
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 428571420
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -428573416
# Row 77
$ws.Range("H77").Value = 20000
$ws.Range("J77").Value = 20000
$ws.Range("L77").Value = 180000
$ws.Range("N77").Value = -190608
# Row 78
$ws.Range("H78").Value = 142857140
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 142857140
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 1285714260
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -1285724244
# Row 82
$ws.Range("H82").Value = 1500
$ws.Range("I82").Value = 1500
$ws.Range("K82").Value = 4500
$ws.Range("M82").Value = -4094
# Row 85
$ws.Range("H85").Value = 1500
$ws.Range("I85").Value = 1500
$ws.Range("K85").Value = 4500
$ws.Range("M85").Value = -3096
# Row 88
$ws.Range("H88").Value = 4087.5
$ws.Range("J88").Value = 4087.5
$ws.Range("L88").Value = 12262.5
$ws.Range("N88").Value = -13118.5
# Row 91
$ws.Range("H91").Value = 4087.5
$ws.Range("J91").Value = 4087.5
$ws.Range("L91").Value = 12262.5
$ws.Range("N91").Value = -15226.5
# Row 105
$ws.Range("H105").Value = 6939.3335
$ws.Range("J105").Value = 6939.3335
$ws.Range("L105").Value = 20818.0005
$ws.Range("N105").Value = -26060.0005
# Row 108
$ws.Range("H108").Value = 100400
$ws.Range("I108").Value = 100400
$ws.Range("K108").Value = 301200
$ws.Range("M108").Value = -298320
# Row 109
$ws.Range("H109").Value = 2314.2727
$ws.Range("I109").Value = 963.5
$ws.Range("J109").Value = 2614.4443
$ws.Range("K109").Value = 2890.5
$ws.Range("L109").Value = 7843.3329
$ws.Range("M109").Value = -1850.5
$ws.Range("N109").Value = -9923.332900000001
# Row 117
$ws.Range("H117").Value = 19613522
$ws.Range("J117").Value = 27785684
$ws.Range("L117").Value = 83357052
$ws.Range("N117").Value = -83363936
# Row 134
$ws.Range("H134").Value = 11156.929
$ws.Range("I134").Value = 13910.444
$ws.Range("J134").Value = 9852.632
$ws.Range("K134").Value = 41731.33199999999
$ws.Range("L134").Value = 29557.896
$ws.Range("M134").Value = -36661.33199999999
$ws.Range("N134").Value = -39697.896

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 29
$ws.Range("H29").Value = 200001600
$ws.Range("I29").Value = 1000000000
$ws.Range("K29").Value = 1000000000
$ws.Range("M29").Value = -999999710
# Row 123
$ws.Range("H123").Value = 27290.062
$ws.Range("J123").Value = 27290.062
$ws.Range("L123").Value = 27290.062
$ws.Range("N123").Value = -32190.062

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 4354.039
$ws.Range("I136").Value = 2307.6667
$ws.Range("J136").Value = 8105.722
$ws.Range("K136").Value = 6923.000100000001
$ws.Range("L136").Value = 24317.166
$ws.Range("M136").Value = -4373.000100000001
$ws.Range("N136").Value = -29417.166

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 47
$ws.Range("H47").Value = 15000
$ws.Range("J47").Value = 15000
$ws.Range("L47").Value = 15000
$ws.Range("N47").Value = -16144
